$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (columns A-R); column S stays empty for all rows.
$data = @(
    @(128, 2731, 130, 1544, 79, 4484, 62, 68, 130, 27973, 19156, 4333, 0, 23640, 0, 0, 1674, 0),
    @(129, 2731, 130, 1544, 79, 4484, 62, 68, 130, 27973, 19156, 4333, 0, 23640, 0, 0, 1674, 0),
    @(130, 2731, 166, 2151, 90, 5138, 75, 91, 166, 30453, 20541, 4774, 654, 25679, 11, 0, 2317, 2039)
)

$startRow = 129
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
    # Column S (19) stays empty for these rows, but still gets a cell
    # placeholder in the sheet (matching the existing rows above it).
    # Touching a no-op border keeps it present without allocating a new style.
    $ws.Cells.Item($row, 19).Borders.LineStyle = -4142
}
